$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.720.53"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "3.463.69"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'578.55"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").Value = "'146.51"
$ws.Range("E6").Value = "  +3.82%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.483"
$ws.Range("E8").Value = "  +2.07%  "
$ws.Range("D9").Value = "'7.63"
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("E10").Value = "  +1.89%  "
$ws.Range("D11").Value = "'0.400"
$ws.Range("E11").Value = "  +3.73%  "
$ws.Range("D12").Value = "4.058.54"
$ws.Range("D13").Value = "'29.87"
$ws.Range("E13").Value = "  +5.30%  "
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").Value = "3.470.35"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "62.838.59"
$ws.Range("E17").Value = "  +1.96%  "
$ws.Range("D18").Value = "'6.36"
$ws.Range("E18").Value = "  +3.74%  "
$ws.Range("D19").Value = "'14.42"
$ws.Range("E19").Value = "  +5.94%  "
$ws.Range("D20").Value = "'9.26"
$ws.Range("E20").Value = "  +2.74%  "
$ws.Range("D21").Value = "'388.68"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("D22").Value = "'0.565"
$ws.Range("E22").Value = "  +2.70%  "
$ws.Range("D23").Value = "'74.99"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "3.609.47"
$ws.Range("E25").Value = "  +2.37%  "
$ws.Range("E26").Value = "  +1.44%  "
$ws.Range("D27").Value = "'0.178"
$ws.Range("E27").Value = "  -8.04%  "
$ws.Range("D28").Value = "'7.62"
$ws.Range("E28").Value = "  +4.81%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "'8.13"
$ws.Range("E30").Value = "  +1.21%  "
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "'1.39"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("D34").Value = "'23.75"
$ws.Range("E34").Value = "  +1.84%  "
$ws.Range("D35").Value = "'7.09"
$ws.Range("E35").Value = "  +2.75%  "
$ws.Range("D36").Value = "'5.28"
$ws.Range("E36").Value = "  +5.00%  "
$ws.Range("B37").Value = "EnergySwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D37").Value = "'31.46"
$ws.Range("E37").Value = "  +20.16%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'1.57"
$ws.Range("E38").Value = "  +6.72%  "
$ws.Range("D39").Value = "'169.73"
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("D40").Value = "3.504.37"
$ws.Range("E40").Value = "  +2.32%  "
$ws.Range("D41").Value = "'0.0770"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("D42").Value = "'0.799"
$ws.Range("E42").Value = "  +2.42%  "
$ws.Range("D43").Value = "'4.49"
$ws.Range("E43").Value = "  +1.73%  "
$ws.Range("D44").Value = "'42.22"
$ws.Range("E45").Value = "  +3.38%  "
$ws.Range("E46").Value = "  +2.70%  "
$ws.Range("D47").Value = "2.596.68"
$ws.Range("E47").Value = "  +4.88%  "
$ws.Range("D48").Value = "'23.41"
$ws.Range("E48").Value = "  +2.60%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'6.78"
$ws.Range("E49").Value = "  +1.80%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "'2.20"
$ws.Range("E50").Value = "  +9.51%  "
$ws.Range("E51").Value = "  +0.03%  "
